$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 387, shifting the
# existing rows 387-403 down to 389-405 (matches the diff, which shows
# the old row-387..403 data reappearing one/two rows further down, plus
# two brand-new rows of data at the top of that block).
$ws.Rows.Item(387).Insert()
$ws.Rows.Item(387).Insert()

# New row 387
$ws.Range("A387").Value = 8
$ws.Range("B387").Value = "Terminal La Palmera de La Serena"
$ws.Range("C387").Value = "Coquimbo"
$ws.Range("D387").Value = 44939
$ws.Range("E387").Value = 4
$ws.Range("F387").Value = 100112003
$ws.Range("G387").Value = "Ajo"
$ws.Range("H387").Value = "Chino"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 400
$ws.Range("K387").Value = 15000
$ws.Range("L387").Value = 16000
$ws.Range("M387").Value = 15500
$ws.Range("N387").Value = "$/caja 10 kilos"
$ws.Range("O387").Value = "China"
$ws.Range("P387").Value = 1550
$ws.Range("Q387").Value = 10
$ws.Range("R387").Value = "Hortaliza"

# New row 388
$ws.Range("A388").Value = 8
$ws.Range("B388").Value = "Terminal La Palmera de La Serena"
$ws.Range("C388").Value = "Coquimbo"
$ws.Range("D388").Value = 44939
$ws.Range("E388").Value = 4
$ws.Range("F388").Value = 100112003
$ws.Range("G388").Value = "Ajo"
$ws.Range("H388").Value = "Chino"
$ws.Range("I388").Value = "Primera"
$ws.Range("J388").Value = 440
$ws.Range("K388").Value = 18000
$ws.Range("L388").Value = 19000
$ws.Range("M388").Value = 18500
$ws.Range("N388").Value = "$/malla 10 kilos"
$ws.Range("O388").Value = "China"
$ws.Range("P388").Value = 1850
$ws.Range("Q388").Value = 10
$ws.Range("R388").Value = "Hortaliza"
